$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# B12: number of investment initialization years: 0 -> 3
$ws.Range("B12").Value = 3

# C12: update description text for investment_initialization_years
$ws.Range("C12").Value = ".Should be between 0 and 3. If this is 3 then the initialization loop should not be executed, if it is 0 then it the initialization loop is executed"

# Update active cell/selection on the Coupling Parameters sheet
$ws.Activate()
$ws.Range("C14").Select() | Out-Null
